# Scene.xlsx - "new implement of scene format"
# The dungeon-entrance rows (13020xxx, rows 23-30) get a distinct TilePath
# ("Q" column) value, prefixed with "dg", separating the dungeon tile path
# from the scene's own Url ("P" column) which keeps the original name.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scene")

$ws.Range("Q23").Value = "dgforestmaze"
$ws.Range("Q24").Value = "dgforestinner"
$ws.Range("Q25").Value = "dgpersepalace1"
$ws.Range("Q26").Value = "dgpersepalace2"
$ws.Range("Q27").Value = "dgpersepalace3"
$ws.Range("Q28").Value = "dgviliage1"
$ws.Range("Q29").Value = "dgviliage2"
$ws.Range("Q30").Value = "dgviliage3"

# Reflect the author's last selection before saving.
$ws.Range("Q24").Select()
